$d = $word.ActiveDocument

$replacements = @(
    @{old = "63×30="; new = "27×82="},
    @{old = "37×85="; new = "74×36="},
    @{old = "95×89="; new = "69×61="},
    @{old = "34×36="; new = "25×70="},
    @{old = "52×34="; new = "26×60="},
    @{old = "51×99="; new = "71×60="},
    @{old = "25×75="; new = "33×36="},
    @{old = "21×99="; new = "97×44="},
    @{old = "37×46="; new = "31×23="},
    @{old = "70×76="; new = "61×55="},
    @{old = "16×42="; new = "44×18="},
    @{old = "94×47="; new = "68×63="},
    @{old = "55×65="; new = "80×46="},
    @{old = "85×41="; new = "98×39="},
    @{old = "81×98="; new = "36×30="},
    @{old = "64×72="; new = "99×75="},
    @{old = "97×93="; new = "12×37="},
    @{old = "33×47="; new = "72×49="},
    @{old = "33×71="; new = "43×78="},
    @{old = "70×63="; new = "81×85="},
    @{old = "43×79="; new = "59×45="},
    @{old = "51×58="; new = "99×63="},
    @{old = "22×99="; new = "45×97="},
    @{old = "86×90="; new = "58×27="},
    @{old = "86×34="; new = "36×77="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
